$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) updates for the refreshed crypto snapshot
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.821.77"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.549.51"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.86"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.66"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.546.39"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  +5.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.148.41"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000198"
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.26"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.561.03"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.630.92"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.28"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.35"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.81"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.688.70"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -5.31%  "
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.554.13"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.59"
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.66"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.61"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0834"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.866"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.48"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -6.31%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.18"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.44"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.926"
$ws.Range("E51").Value = "  -2.95%  "
